$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 11 rows of label/value pairs in columns A and B.
# We only want to keep the values (column B) and lay them out across a
# single row, columns A through K.
$v1 = $ws.Range("B1").Value2
$v2 = $ws.Range("B2").Value2
$v3 = $ws.Range("B3").Value2
$v4 = $ws.Range("B4").Value2
$v5 = $ws.Range("B5").Value2
$v6 = $ws.Range("B6").Value2
$v7 = $ws.Range("B7").Value2
$v8 = $ws.Range("B8").Value2
$v9 = $ws.Range("B9").Value2
$v10 = $ws.Range("B10").Value2
$v11 = $ws.Range("B11").Value2

# Clear out the old 2-column x 11-row block entirely.
$ws.Range("A1:B11").Clear()

# Write the values back out as a single row, columns A through K.
$ws.Range("A1").Value = $v1
$ws.Range("B1").Value = $v2
$ws.Range("C1").Value = $v3
$ws.Range("D1").Value = $v4
$ws.Range("E1").Value = $v5
$ws.Range("F1").Value = $v6
$ws.Range("G1").Value = $v7
$ws.Range("H1").Value = $v8
$ws.Range("I1").Value = $v9
$ws.Range("J1").Value = $v10
$ws.Range("K1").Value = $v11
